# Updated symbol list on Sat Dec 31 05:33:12 UTC 2022 with GitHub Actions
#
# Applies the per-cell value updates captured in the diff against
# cryptos.xlsx (Sheet1). Column D ("Price") holds numeric-looking values
# that are stored as text in the workbook, so those are written with a
# leading apostrophe (quote-prefix) to keep them text instead of letting
# Excel auto-convert them to numbers. Columns B/C/E hold ordinary text and
# are set directly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column D price updates (kept as text via quote-prefix) ---
$ws.Range("D2").Formula  = "'245.62"
$ws.Range("D3").Formula  = "'25.47"
$ws.Range("D4").Formula  = "'5.090"
$ws.Range("D5").Formula  = "'0.05572"
$ws.Range("D6").Formula  = "'6.497"
$ws.Range("D7").Formula  = "'3.019"
$ws.Range("D8").Formula  = "'0.8190"
$ws.Range("D9").Formula  = "'0.8441"
$ws.Range("D10").Formula = "'0.1342"
$ws.Range("D11").Formula = "'0.06950"
$ws.Range("D12").Formula = "'0.03188"
$ws.Range("D13").Formula = "'0.02869"
$ws.Range("D14").Formula = "'0.09379"
$ws.Range("D15").Formula = "'0.001521"
$ws.Range("D16").Formula = "'0.0005961"
$ws.Range("D17").Formula = "'0.006077"
$ws.Range("D18").Formula = "'3.499"
$ws.Range("D22").Formula = "'3.768"
$ws.Range("D26").Formula = "'0.004635"
$ws.Range("D27").Formula = "'0.00009705"

# --- Row 16 / 27 "Volume(1h)" label tweaks ---
$ws.Range("E16").Value = "15OneONEWorstin24h"
$ws.Range("E27").Value = "26NitroExNTXBestin24h"

# --- Rows 41/42 swapped coin entries (KickToken <-> BKEXToken) ---
$ws.Range("B41").Value   = "BKEXToken"
$ws.Range("C41").Value   = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D41").Formula = "'0.1052"
$ws.Range("E41").Value   = "40BKEXTokenBKK"

$ws.Range("B42").Value   = "KickToken"
$ws.Range("C42").Value   = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D42").Formula = "'0.006214"
$ws.Range("E42").Value   = "41KickTokenKICK"

# --- Remaining column D price updates ---
$ws.Range("D43").Formula = "'0.002001"
$ws.Range("D44").Formula = "'0.008319"
$ws.Range("D45").Formula = "'0.00005307"

# --- Row 47 label tweak ---
$ws.Range("E47").Value = "46CoinbaseStockTokenCOIN"

# --- Final column D price updates ---
$ws.Range("D48").Formula = "'0.002123"
$ws.Range("D49").Formula = "'0.00002101"
$ws.Range("D50").Formula = "'0.0002001"
